# Insert two new price records right before the existing row 375
# ("Feria Lagunitas de Puerto Montt" / Zanahoria data). This shifts every
# subsequent row down by two, which is what the target diff shows (the
# sheet's used range grows from A1:R461 to A1:R463, and the two rows that
# used to be the last ones now reappear as the new last two rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 375, pushing old row 375.. down to 377..
$ws.Rows("375:376").Insert()

# --- New row 375 ---
$ws.Range("A375").Value = 4
$ws.Range("B375").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C375").Value = "Los Lagos"
$ws.Range("D375").Value = 44889
$ws.Range("E375").Value = 10
$ws.Range("F375").Value = 100114013
$ws.Range("G375").Value = "Zanahoria"
$ws.Range("H375").Value = "Sin especificar"
$ws.Range("I375").Value = "Primera"
$ws.Range("J375").Value = 200
$ws.Range("K375").Value = 15000
$ws.Range("L375").Value = 15000
$ws.Range("M375").Value = 15000
$ws.Range("N375").Value = "$/saco 20 kilos"
$ws.Range("O375").Value = "Región Metropolitana"
$ws.Range("P375").Value = 750
$ws.Range("Q375").Value = 20
$ws.Range("R375").Value = "Hortaliza"

# --- New row 376 ---
$ws.Range("A376").Value = 4
$ws.Range("B376").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C376").Value = "Los Lagos"
$ws.Range("D376").Value = 44889
$ws.Range("E376").Value = 10
$ws.Range("F376").Value = 100114013
$ws.Range("G376").Value = "Zanahoria"
$ws.Range("H376").Value = "Sin especificar"
$ws.Range("I376").Value = "Segunda"
$ws.Range("J376").Value = 200
$ws.Range("K376").Value = 12000
$ws.Range("L376").Value = 12000
$ws.Range("M376").Value = 12000
$ws.Range("N376").Value = "$/saco 20 kilos"
$ws.Range("O376").Value = "Región Metropolitana"
$ws.Range("P376").Value = 600
$ws.Range("Q376").Value = 20
$ws.Range("R376").Value = "Hortaliza"
